{"js": "// Update the \"C\u00e1c b\u01b0\u1edbc th\u1ef1c hi\u1ec7n\" procedure list:\n//  - B\u01b0\u1edbc 3: drop the \"(histogram, KDE plot)\" qualifier.\n//  - B\u01b0\u1edbc 4: generalize \"GRE, TOEFL, CGPA\" to \"c\u00e1c thu\u1ed9c t\u00ednh\".\n//  - Remove the old B\u01b0\u1edbc 5 (heatmap) and B\u01b0\u1edbc 6 (boxplot) steps entirely.\n//  - The old B\u01b0\u1edbc 7 step is removed too, but its slot/paragraph is reused:\n//    it becomes the new \"B\u01b0\u1edbc 5\" step carrying the heatmap text that used\n//    to belong to the removed step 5.\nconst body = context.document.body;\n\n// --- Step 3: trim the parenthetical from the input-distribution bullet ---\nconst step3Matches = body.search(\n  \"Tr\u1ef1c quan h\u00f3a ph\u00e2n ph\u1ed1i c\u00e1c bi\u1ebfn \u0111\u1ea7u v\u00e0o (histogram, KDE plot).\",\n  { matchCase: true }\n);\nstep3Matches.load(\"items\");\nawait context.sync();\nif (step3Matches.items.length === 0) {\n  throw new Error(\"Could not find B\u01b0\u1edbc 3 text to update\");\n}\nstep3Matches.items[0].insertText(\n  \"Tr\u1ef1c quan h\u00f3a ph\u00e2n ph\u1ed1i c\u00e1c bi\u1ebfn \u0111\u1ea7u v\u00e0o .\",\n  Word.InsertLocation.replace\n);\n\n// --- Step 4: generalize the scatter-plot variable list ---\nconst step4Matches = body.search(\n  \"V\u1ebd scatter plot gi\u1eefa GRE, TOEFL, CGPA v\u1edbi Chance of Admit.\",\n  { matchCase: true }\n);\nstep4Matches.load(\"items\");\nawait context.sync();\nif (step4Matches.items.length === 0) {\n  throw new Error(\"Could not find B\u01b0\u1edbc 4 text to update\");\n}\nstep4Matches.items[0].insertText(\n  \"V\u1ebd scatter plot gi\u1eefa c\u00e1c thu\u1ed9c t\u00ednh v\u1edbi Chance of Admit.\",\n  Word.InsertLocation.replace\n);\n\n// --- Remove the old B\u01b0\u1edbc 5 (heatmap) and B\u01b0\u1edbc 6 (boxplot) paragraphs ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet oldStep5Paragraph = null;\nlet oldStep6Paragraph = null;\nlet oldStep7Paragraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Ph\u00e2n t\u00edch t\u01b0\u01a1ng quan b\u1eb1ng bi\u1ec3u \u0111\u1ed3 heatmap.\") !== -1) {\n    oldStep5Paragraph = paragraphs.items[i];\n  } else if (\n    text.indexOf(\n      \"So s\u00e1nh kh\u1ea3 n\u0103ng tr\u00fang tuy\u1ec3n gi\u1eefa nh\u00f3m c\u00f3 v\u00e0 kh\u00f4ng c\u00f3 kinh nghi\u1ec7m nghi\u00ean c\u1ee9u (boxplot).\"\n    ) !== -1\n  ) {\n    oldStep6Paragraph = paragraphs.items[i];\n  } else if (\n    text.indexOf(\"\u0110\u00e1nh gi\u00e1 \u1ea3nh h\u01b0\u1edfng c\u1ee7a University Rating, SOP v\u00e0 LOR.\") !== -1\n  ) {\n    oldStep7Paragraph = paragraphs.items[i];\n  }\n}\nif (!oldStep5Paragraph || !oldStep6Paragraph || !oldStep7Paragraph) {\n  throw new Error(\"Could not locate B\u01b0\u1edbc 5/6/7 paragraphs\");\n}\n\noldStep5Paragraph.delete();\noldStep6Paragraph.delete();\nawait context.sync();\n\n// --- Turn the former B\u01b0\u1edbc 7 paragraph into the new B\u01b0\u1edbc 5 ---\nconst labelMatches = oldStep7Paragraph.search(\"B\u01b0\u1edbc 7: \", { matchCase: true });\nlabelMatches.load(\"items\");\nawait context.sync();\nif (labelMatches.items.length === 0) {\n  throw new Error(\"Could not find 'B\u01b0\u1edbc 7: ' label to rename\");\n}\nlabelMatches.items[0].insertText(\"B\u01b0\u1edbc 5:\", Word.InsertLocation.replace);\n\nconst bodyMatches = oldStep7Paragraph.search(\n  \"\u0110\u00e1nh gi\u00e1 \u1ea3nh h\u01b0\u1edfng c\u1ee7a University Rating, SOP v\u00e0 LOR.\",\n  { matchCase: true }\n);\nbodyMatches.load(\"items\");\nawait context.sync();\nif (bodyMatches.items.length === 0) {\n  throw new Error(\"Could not find B\u01b0\u1edbc 7 body text to replace\");\n}\nbodyMatches.items[0].insertText(\n  \" Ph\u00e2n t\u00edch t\u01b0\u01a1ng quan b\u1eb1ng bi\u1ec3u \u0111\u1ed3 heatmap.\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Update the \"C\u00e1c b\u01b0\u1edbc th\u1ef1c hi\u1ec7n\" procedure list:\n#  - B\u01b0\u1edbc 3: drop the \"(histogram, KDE plot)\" qualifier.\n#  - B\u01b0\u1edbc 4: generalize \"GRE, TOEFL, CGPA\" to \"c\u00e1c thu\u1ed9c t\u00ednh\".\n#  - Remove the old B\u01b0\u1edbc 5 (heatmap) and B\u01b0\u1edbc 6 (boxplot) steps entirely.\n#  - The old B\u01b0\u1edbc 7 step is removed too, but its paragraph is reused: it\n#    becomes the new \"B\u01b0\u1edbc 5\" step carrying the heatmap text that used to\n#    belong to the removed step 5.\n\n$d = $word.ActiveDocument\n\n# --- Step 3: trim the parenthetical from the input-distribution bullet ---\n$oldStep3 = \"Tr\u1ef1c quan h\u00f3a ph\u00e2n ph\u1ed1i c\u00e1c bi\u1ebfn \u0111\u1ea7u v\u00e0o (histogram, KDE plot).\"\n$newStep3 = \"Tr\u1ef1c quan h\u00f3a ph\u00e2n ph\u1ed1i c\u00e1c bi\u1ebfn \u0111\u1ea7u v\u00e0o .\"\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Text = $oldStep3\n$found3 = $find3.Execute($oldStep3, $false, $false, $false, $false, $false, $true, 1, $false, $newStep3, 2)\nif (-not $found3) {\n    throw \"Could not find B\u01b0\u1edbc 3 text to update\"\n}\n\n# --- Step 4: generalize the scatter-plot variable list ---\n$oldStep4 = \"V\u1ebd scatter plot gi\u1eefa GRE, TOEFL, CGPA v\u1edbi Chance of Admit.\"\n$newStep4 = \"V\u1ebd scatter plot gi\u1eefa c\u00e1c thu\u1ed9c t\u00ednh v\u1edbi Chance of Admit.\"\n$find4 = $d.Content.Find\n$find4.ClearFormatting()\n$find4.Text = $oldStep4\n$found4 = $find4.Execute($oldStep4, $false, $false, $false, $false, $false, $true, 1, $false, $newStep4, 2)\nif (-not $found4) {\n    throw \"Could not find B\u01b0\u1edbc 4 text to update\"\n}\n\n# --- Remove the old B\u01b0\u1edbc 5 (heatmap) and B\u01b0\u1edbc 6 (boxplot) paragraphs ---\n$step5Marker = \"Ph\u00e2n t\u00edch t\u01b0\u01a1ng quan b\u1eb1ng bi\u1ec3u \u0111\u1ed3 heatmap.\"\n$step6Marker = \"So s\u00e1nh kh\u1ea3 n\u0103ng tr\u00fang tuy\u1ec3n gi\u1eefa nh\u00f3m c\u00f3 v\u00e0 kh\u00f4ng c\u00f3 kinh nghi\u1ec7m nghi\u00ean c\u1ee9u (boxplot).\"\n$step7Marker = \"\u0110\u00e1nh gi\u00e1 \u1ea3nh h\u01b0\u1edfng c\u1ee7a University Rating, SOP v\u00e0 LOR.\"\n\n$step5Paragraph = $null\n$step6Paragraph = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.Contains($step5Marker)) {\n        $step5Paragraph = $p\n    } elseif ($t.Contains($step6Marker)) {\n        $step6Paragraph = $p\n    }\n}\nif ($null -eq $step5Paragraph -or $null -eq $step6Paragraph) {\n    throw \"Could not locate B\u01b0\u1edbc 5/6 paragraphs\"\n}\n\n# Delete the later paragraph first so the earlier paragraph's position\n# stays valid for the second delete.\n$step6Paragraph.Range.Delete()\n$step5Paragraph.Range.Delete()\n\n# Re-locate the (former) B\u01b0\u1edbc 7 paragraph now that it has shifted up;\n# ranges captured before the deletes above do not track the move.\n$step7Paragraph = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.Contains($step7Marker)) {\n        $step7Paragraph = $p\n    }\n}\nif ($null -eq $step7Paragraph) {\n    throw \"Could not locate the (former) B\u01b0\u1edbc 7 paragraph after deletion\"\n}\n$step7Range = $step7Paragraph.Range\n\n# --- Turn the former B\u01b0\u1edbc 7 paragraph into the new B\u01b0\u1edbc 5 ---\n$findLabel = $step7Range.Find\n$findLabel.ClearFormatting()\n$findLabel.Text = \"B\u01b0\u1edbc 7: \"\n$foundLabel = $findLabel.Execute(\"B\u01b0\u1edbc 7: \", $false, $false, $false, $false, $false, $true, 1, $false, \"B\u01b0\u1edbc 5:\", 2)\nif (-not $foundLabel) {\n    throw \"Could not find 'B\u01b0\u1edbc 7: ' label to rename\"\n}\n\n$findBody = $step7Range.Find\n$findBody.ClearFormatting()\n$findBody.Text = $step7Marker\n$foundBody = $findBody.Execute($step7Marker, $false, $false, $false, $false, $false, $true, 1, $false, \" Ph\u00e2n t\u00edch t\u01b0\u01a1ng quan b\u1eb1ng bi\u1ec3u \u0111\u1ed3 heatmap.\", 2)\nif (-not $foundBody) {\n    throw \"Could not find B\u01b0\u1edbc 7 body text to replace\"\n}\n"}
